$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values look like numbers to Excel's auto-detect
# (e.g. "21.22", "0.4290") even though the source column stores them as plain
# text (see e.g. the untouched "1.007" cells already in the sheet). Force those
# specific cells to Text before writing so the literal string is preserved,
# then restore the default "Normal" style so no stray number-format/quote-prefix
# styling is left behind.
$textCells = @(
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D12",
    "D13",
    "D16",
    "D17",
    "D18",
    "D20",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.906.57"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "1.810.37"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "312.63"
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "0.4290"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "0.3691"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").Value = "0.07264"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("D10").Value = "0.8635"
$ws.Range("E10").Value = "  +4.45%  "
$ws.Range("D11").Value = "2.048.99"
$ws.Range("E11").Value = "  +17.89%  "
$ws.Range("D12").Value = "21.22"
$ws.Range("E12").Value = "  +5.93%  "
$ws.Range("D13").Value = "6.626"
$ws.Range("E13").Value = "  +5.08%  "
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "80.75"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "0.000008914"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "15.17"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").Value = "26.962.20"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").Value = "5.194"
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "2.285.94"
$ws.Range("E24").Value = "  +16.55%  "
$ws.Range("D25").Value = "154.05"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "1.884"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "18.35"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").Value = "5.232"
$ws.Range("E28").Value = "  +4.75%  "
$ws.Range("D29").Value = "1.911"
$ws.Range("E29").Value = "  +16.48%  "
$ws.Range("D30").Value = "114.76"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").Value = "0.7395"
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("D33").Value = "1.157"
$ws.Range("E33").Value = "  +6.54%  "
$ws.Range("D34").Value = "4.426"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "1.123"
$ws.Range("E37").Value = "  +5.76%  "
$ws.Range("D38").Value = "0.05228"
$ws.Range("E38").Value = "  +3.30%  "
$ws.Range("D39").Value = "0.01923"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("D40").Value = "0.5080"
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("D41").Value = "2.754"
$ws.Range("E41").Value = "  +12.00%  "
$ws.Range("D42").Value = "0.1647"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").Value = "6.451"
$ws.Range("E43").Value = "  +5.61%  "
$ws.Range("D44").Value = "8.295"
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("D45").Value = "107.23"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").Value = "10.36"
$ws.Range("E46").Value = "  +4.39%  "
$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "1.647"
$ws.Range("E48").Value = "  +5.75%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.4558"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06277"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").Value = "1.809"
$ws.Range("E51").Value = "  +6.47%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
